$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-looking Price (D-column) cells to stay as text, matching the
# original inline-string cell type, then restore the default "Normal" style so
# no stray formatting is introduced.

$ws.Range("D2").Value = "60.259.95"
$ws.Range("E2").Value = "  +4.09%  "
$ws.Range("D3").Value = "2.430.17"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.53%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("E9").Value = "  +4.55%  "
$ws.Range("E10").Value = "  +4.20%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.38%  "
$ws.Range("D14").Value = "2.863.10"
$ws.Range("E14").Value = "  +3.23%  "
$ws.Range("D15").Value = "60.183.87"
$ws.Range("E15").Value = "  +4.07%  "
$ws.Range("D17").Value = "2.392.52"
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.15%  "
$ws.Range("E19").Value = "  +3.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "333.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.33%  "
$ws.Range("E24").Value = "  +3.22%  "
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "0.0₃0789"
$ws.Range("E28").Value = "  +7.11%  "
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("E32").Value = "  +2.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +5.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "325.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.47%  "
$ws.Range("E40").Value = "  +11.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "39.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "140.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0526"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.61%  "
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.415"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.86%  "
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("E51").Value = "  -0.17%  "
